{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph that separates it from the \"Requisitos\" text\n// above it) from the end of the document body.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nconst VER_NO_JUPITER = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst COPYRIGHT_PREFIX = \"\\u00A9 2020 . Contact: luizeleno@usp.br.\";\n\n// Locate the \"Ver no Jupiter ...\" paragraph; the \"\u00a9 2020 ...\" paragraph is\n// expected to immediately follow it, and a single blank paragraph\n// immediately precedes it (separating it from the preceding content).\nlet verIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === VER_NO_JUPITER) {\n    verIndex = i;\n    break;\n  }\n}\n\nif (verIndex !== -1) {\n  const toDelete = [];\n\n  // The paragraph with \"Ver no Jupiter ...\"\n  toDelete.push(items[verIndex]);\n\n  // The following paragraph, if it is the copyright/footer line.\n  if (\n    verIndex + 1 < items.length &&\n    items[verIndex + 1].text.indexOf(COPYRIGHT_PREFIX) === 0\n  ) {\n    toDelete.push(items[verIndex + 1]);\n  }\n\n  // The blank separator paragraph right before \"Ver no Jupiter ...\".\n  if (verIndex - 1 >= 0 && items[verIndex - 1].text === \"\") {\n    toDelete.push(items[verIndex - 1]);\n  }\n\n  // Delete in reverse document order so earlier deletions don't shift the\n  // position of paragraphs we still need to remove.\n  toDelete.sort((a, b) => items.indexOf(b) - items.indexOf(a));\n  for (const p of toDelete) {\n    p.delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n# (and the blank paragraph that separates it from the \"Requisitos\" text\n# above it) from the end of the document body.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphIndexAt($doc, $pos) {\n    $count = $doc.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        if ($p.Range.Start -le $pos -and $pos -lt $p.Range.End) {\n            return $i\n        }\n    }\n    return -1\n}\n\n$searchRange = $d.Content\n$searchRange.Find.Text = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$found = $searchRange.Find.Execute()\n\nif ($found) {\n    $verIdx = Get-ParagraphIndexAt $d $searchRange.Start\n\n    $startIdx = $verIdx\n    $endIdx = $verIdx\n\n    # Absorb the blank separator paragraph immediately before, if present.\n    if ($verIdx - 1 -ge 1) {\n        $prev = $d.Paragraphs.Item($verIdx - 1)\n        if ($prev.Range.Text.Trim() -eq \"\") {\n            $startIdx = $verIdx - 1\n        }\n    }\n\n    # Absorb the copyright/footer paragraph immediately after, if present.\n    if ($verIdx + 1 -le $d.Paragraphs.Count) {\n        $nxt = $d.Paragraphs.Item($verIdx + 1)\n        $copyPrefix = \"\u00a9 2020\"\n        if ($nxt.Range.Text.StartsWith($copyPrefix)) {\n            $endIdx = $verIdx + 1\n        }\n    }\n\n    $rStart = $d.Paragraphs.Item($startIdx).Range.Start\n    $rEnd = $d.Paragraphs.Item($endIdx).Range.End\n    $delRange = $d.Range($rStart, $rEnd)\n    $delRange.Delete()\n}\n"}
